$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$shp = $hdr.Range.InlineShapes.Item(1)
$conv = $shp.ConvertToShape()
Write-Output ("Title: [" + $conv.Title + "]")
Write-Output ("AlternativeText: [" + $conv.AlternativeText + "]")
try { Write-Output ("ID: " + $conv.ID) } catch { Write-Output "no ID" }
